$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 3267  # was 3266, cell F4
$ws.Cells.Item(11, 6).Value = 19  # was 18, cell F11
$ws.Cells.Item(19, 6).Value = 1538  # was 1531, cell F19
$ws.Cells.Item(20, 6).Value = 328  # was 327, cell F20
$ws.Cells.Item(21, 6).Value = 6043  # was 6041, cell F21
$ws.Cells.Item(22, 6).Value = 35  # was 34, cell F22
$ws.Cells.Item(23, 6).Value = 2341  # was 2340, cell F23
$ws.Cells.Item(27, 6).Value = 4193  # was 4189, cell F27
$ws.Cells.Item(28, 6).Value = 3838  # was 3836, cell F28
$ws.Cells.Item(29, 6).Value = 284  # was 283, cell F29
$ws.Cells.Item(30, 6).Value = 89  # was 87, cell F30
$ws.Cells.Item(34, 6).Value = 1021  # was 1020, cell F34
$ws.Cells.Item(36, 6).Value = 82  # was 81, cell F36
$ws.Cells.Item(39, 6).Value = 9  # was 8, cell F39
$ws.Cells.Item(41, 6).Value = 22  # was 21, cell F41
$ws.Cells.Item(42, 6).Value = 578  # was 577, cell F42
$ws.Cells.Item(43, 6).Value = 377  # was 375, cell F43
$ws.Cells.Item(44, 6).Value = 311  # was 310, cell F44
$ws.Cells.Item(45, 6).Value = 1082  # was 1080, cell F45
$ws.Cells.Item(47, 6).Value = 3090  # was 3037, cell F47
$ws.Cells.Item(48, 6).Value = 72  # was 69, cell F48
$ws.Cells.Item(49, 6).Value = 335  # was 331, cell F49
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 6).Value = 115  # was 114, cell F6
$ws.Cells.Item(6, 7).Value = 300  # was 260, cell G6
$ws.Cells.Item(10, 6).Value = 602  # was 601, cell F10
$ws.Cells.Item(14, 6).Value = 79  # was 78, cell F14
$ws.Cells.Item(15, 6).Value = 240  # was 239, cell F15
$ws.Cells.Item(24, 6).Value = 2  # was 1, cell F24
$ws.Cells.Item(28, 6).Value = 5288  # was 5211, cell F28
$ws.Cells.Item(29, 6).Value = 5288  # was 5211, cell F29
$ws.Cells.Item(30, 6).Value = 42  # was 41, cell F30
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(6, 6).Value = 1952  # was 1950, cell F6
$ws.Cells.Item(8, 6).Value = 3028  # was 3027, cell F8
$ws.Cells.Item(13, 6).Value = 2047  # was 2045, cell F13
$ws.Cells.Item(14, 6).Value = 8745  # was 8743, cell F14
$ws.Cells.Item(15, 6).Value = 889  # was 887, cell F15
$ws.Cells.Item(16, 6).Value = 51  # was 50, cell F16
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 1952  # was 1950, cell F4
$ws.Cells.Item(6, 6).Value = 3028  # was 3027, cell F6
$ws.Cells.Item(12, 6).Value = 115  # was 114, cell F12
$ws.Cells.Item(12, 7).Value = 300  # was 260, cell G12
$ws.Cells.Item(13, 6).Value = 19  # was 18, cell F13
$ws.Cells.Item(16, 6).Value = 889  # was 887, cell F16
$ws.Cells.Item(17, 6).Value = 602  # was 601, cell F17
$ws.Cells.Item(18, 6).Value = 51  # was 50, cell F18
$ws.Cells.Item(25, 6).Value = 240  # was 239, cell F25
$ws.Cells.Item(26, 6).Value = 328  # was 327, cell F26
$ws.Cells.Item(27, 6).Value = 6043  # was 6041, cell F27
$ws.Cells.Item(28, 6).Value = 2341  # was 2340, cell F28
$ws.Cells.Item(31, 6).Value = 4193  # was 4189, cell F31
$ws.Cells.Item(32, 6).Value = 284  # was 283, cell F32
$ws.Cells.Item(36, 6).Value = 1021  # was 1020, cell F36
$ws.Cells.Item(41, 6).Value = 2  # was 1, cell F41
$ws.Cells.Item(42, 6).Value = 578  # was 577, cell F42
$ws.Cells.Item(43, 6).Value = 377  # was 375, cell F43
$ws.Cells.Item(44, 6).Value = 311  # was 310, cell F44
$ws.Cells.Item(47, 6).Value = 3090  # was 3038, cell F47
$ws.Cells.Item(48, 6).Value = 72  # was 69, cell F48
$ws.Cells.Item(49, 6).Value = 5288  # was 5211, cell F49
